# Update row 9 (Ano 2025) figures in the faturamento anual sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B9").Value = 3747024.58
$ws.Range("C9").Value = 587360.24
$ws.Range("D9").Value = 4334384.82
$ws.Range("E9").Value = 13.55117887294557
$ws.Range("F9").Value = 86.44882112705442
$ws.Range("G9").Value = -43.23456151717517
$ws.Range("H9").Value = -32.3339153710111
$ws.Range("I9").Value = 37741
$ws.Range("J9").Value = 1603
$ws.Range("K9").Value = 39344
$ws.Range("L9").Value = 27198
$ws.Range("M9").Value = 159.364101036841
$ws.Range("N9").Value = 8.801117662738879
